# Commit: Sat, Apr 11, 2020  5:04:47 AM
#
# The authored diff changes the table-style reference used by the
# cash-flow comparison table on slide 16 from the deck's custom
# "Table_0" style ({C0292CBE-AC90-46CA-AA31-A41D889C54AF}, still defined
# in ppt/tableStyles.xml) to a PowerPoint built-in table-style GUID
# ({ABFB7136-72AB-487E-BCCF-67178BC0BBB0}).
#
# Table styles can't be changed by assigning Table.Style directly (the
# host reports: "Table styles cannot be assigned through a property -
# call Table.ApplyStyle("{GUID}") instead"), so ApplyStyle is used below.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)            # graphicFrame holding the <a:tbl>
$tbl = $sh.Table

$tbl.ApplyStyle("{ABFB7136-72AB-487E-BCCF-67178BC0BBB0}")
